# Update Excel download: remove the "AbteilungKürzel" column (column B).
# The remaining columns ("Titel", "Autor:innen") shift left to become B and C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "AbteilungKürzel" column (column B); this shifts
# the "Titel" and "Autor:innen" columns one position to the left,
# carrying their values and auto-fit column widths with them.
$ws.Columns("B").Delete()

# Update the selection to match the new state captured in the workbook.
$ws.Range("M7").Select() | Out-Null

$wb.Save() | Out-Null
